# Refresh crypto price/volume data (and re-rank rows 26-27) to match the
# latest GitHub Actions snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '42.786.81'
$ws.Cells.Item(2, 5).Value = '  +0.34%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.314.01'
$ws.Cells.Item(3, 5).Value = '  +0.74%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.17%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '311.79'
$ws.Cells.Item(5, 5).Value = '  -1.25%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '106.40'
$ws.Cells.Item(6, 5).Value = '  +3.42%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.14%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.09%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.66%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '40.32'
$ws.Cells.Item(10, 5).Value = '  +2.51%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +0.98%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '8.35'
$ws.Cells.Item(12, 5).Value = '  -1.80%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -0.25%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.992'
$ws.Cells.Item(14, 5).Value = '  -0.76%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '15.31'
$ws.Cells.Item(15, 5).Value = '  -0.04%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '2.664.06'
$ws.Cells.Item(16, 5).Value = '  +0.57%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '2.317.83'
$ws.Cells.Item(17, 5).Value = '  +0.80%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '42.773.68'
$ws.Cells.Item(18, 5).Value = '  +0.41%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '7.48'
$ws.Cells.Item(19, 5).Value = '  -0.76%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.0000106'
$ws.Cells.Item(20, 5).Value = '  -0.45%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '13.43'
$ws.Cells.Item(21, 5).Value = '  -3.56%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '73.71'
$ws.Cells.Item(22, 5).Value = '  -0.36%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '3.50'
$ws.Cells.Item(23, 5).Value = '  -0.97%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '266.92'
$ws.Cells.Item(24, 5).Value = '  +0.19%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.25'
$ws.Cells.Item(25, 5).Value = '  +0.39%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'Filecoin'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '7.84'
$ws.Cells.Item(26, 5).Value = '  +17.86%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'Dai'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '1.01'
$ws.Cells.Item(27, 5).Value = '  -0.17%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '10.97'
$ws.Cells.Item(28, 5).Value = '  +0.71%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +1.44%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '38.74'
$ws.Cells.Item(30, 5).Value = '  +4.56%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '22.40'
$ws.Cells.Item(31, 5).Value = '  -0.78%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '165.73'
$ws.Cells.Item(32, 5).Value = '  +0.21%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.0871'
$ws.Cells.Item(33, 5).Value = '  -1.04%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '2.75'
$ws.Cells.Item(34, 5).Value = '  +6.19%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -0.44%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +2.38%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.112'
$ws.Cells.Item(37, 5).Value = '  -0.62%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +1.39%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.79'
$ws.Cells.Item(39, 5).Value = '  +3.32%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.64'
$ws.Cells.Item(40, 5).Value = '  -2.17%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.20%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '104.66'
$ws.Cells.Item(42, 5).Value = '  +9.24%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '71.15'
$ws.Cells.Item(43, 5).Value = '  +0.78%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.231'
$ws.Cells.Item(44, 5).Value = '  +1.19%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -0.38%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '12.29'
$ws.Cells.Item(46, 5).Value = '  -0.75%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '112.48'
$ws.Cells.Item(47, 5).Value = '  -4.01%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '1.694.84'
$ws.Cells.Item(48, 5).Value = '  +1.93%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '76.52'
$ws.Cells.Item(49, 5).Value = '  -4.21%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +0.47%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -1.07%  '
